$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("B2").Value = "Down (primero"
$ws.Range("C2").Value = "segundo"
$ws.Range("D2").Value = "tercero"

# Delete row 3 entirely (shift cells up)
$ws.Range("A3:D3").Delete()
